# Weekly update: insert two new "Betarraga" price rows (week of D=44855)
# right above the existing price-history block, pushing the old rows
# (291-368) down by two (to 293-370).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 291.
$ws.Rows("291:292").Insert()

# --- New row 291: "Primera" quality, new weekly observation ---
$ws.Range("A291").Value = 4
$ws.Range("B291").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C291").Value = "Los Lagos"
$ws.Range("D291").Value = 44855
$ws.Range("E291").Value = 10
$ws.Range("F291").Value = 100114014
$ws.Range("G291").Value = "Betarraga"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 500
$ws.Range("K291").Value = 1500
$ws.Range("L291").Value = 1500
$ws.Range("M291").Value = 1500
$ws.Range("N291").Value = "`$/paquete 5 unidades"
$ws.Range("O291").Value = "Región del Maule"
$ws.Range("P291").Value = 300
$ws.Range("Q291").Value = 5
$ws.Range("R291").Value = "Hortaliza"

# --- New row 292: "Segunda" quality, same weekly observation ---
$ws.Range("A292").Value = 4
$ws.Range("B292").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C292").Value = "Los Lagos"
$ws.Range("D292").Value = 44855
$ws.Range("E292").Value = 10
$ws.Range("F292").Value = 100114014
$ws.Range("G292").Value = "Betarraga"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Segunda"
$ws.Range("J292").Value = 500
$ws.Range("K292").Value = 1200
$ws.Range("L292").Value = 1200
$ws.Range("M292").Value = 1200
$ws.Range("N292").Value = "`$/paquete 5 unidades"
$ws.Range("O292").Value = "Región del Maule"
$ws.Range("P292").Value = 240
$ws.Range("Q292").Value = 5
$ws.Range("R292").Value = "Hortaliza"
